$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.563.60'
$ws.Range('E2').Value = '  +1.86%  '
$ws.Range('D3').Value = '2.036.94'
$ws.Range('E3').Value = '  +1.45%  '
$ws.Range('E4').Value = '  -0.73%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.21%  '
$ws.Range('E6').Value = '  +1.48%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.18'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.09%  '
$ws.Range('E9').Value = '  +1.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0803'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.07%  '
$ws.Range('E11').Value = '  -0.59%  '
$ws.Range('D12').Value = '2.337.19'
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.41'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.29'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.31%  '
$ws.Range('E15').Value = '  +2.50%  '
$ws.Range('E16').Value = '  +0.97%  '
$ws.Range('D17').Value = '2.027.80'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '37.484.04'
$ws.Range('E18').Value = '  +1.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.03'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Value = '0.0₃0827'
$ws.Range('E21').Value = '  +1.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '223.14'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.83%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.25'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.40'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('E28').Value = '  +6.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.79'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.19%  '
$ws.Range('E30').Value = '  +0.38%  '
$ws.Range('E31').Value = '  +1.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.50'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0607'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.01%  '
$ws.Range('E34').Value = '  +2.08%  '
$ws.Range('E35').Value = '  +8.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.33'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.26'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.55%  '
$ws.Range('B38').Value = 'THORChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.78'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.73%  '
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').Value = '1.480.39'
$ws.Range('E40').Value = '  +0.21%  '
$ws.Range('E41').Value = '  -0.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0939'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.00%  '
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.84'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.50%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '95.31'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.59'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.95%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.20'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +17.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.11'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.62%  '
$ws.Range('E48').Value = '  +1.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.12'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.51%  '
$ws.Range('E50').Value = '  +1.39%  '
$ws.Range('D51').Value = '2.224.81'
$ws.Range('E51').Value = '  +1.02%  '
